$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C:L for rows 2-25
$rows = @(
    @(0.2424750669258486, 0.2064003509529613, 0.1663423400497983, 1.315347342557217, 0.7166479277978084, 0.8322969081920206, 0.9077249771838432, 0.1802786961482212, 1.984705554895072, 0.2151901350040646),
    @(0.2375028024586499, 0.2013032869760707, 0.1648570631337982, 1.328412596869562, 0.727730312184633, 0.842407897739001, 0.9174652689456764, 0.1804069674720736, 1.766695578561496, 0.2145916048839638),
    @(0.2345216179930105, 0.1982205667129193, 0.1640004749258388, 1.337260014341915, 0.7351228804764744, 0.8490510191961889, 0.9240295970885768, 0.1805585059460455, 1.632334867622376, 0.2142992589824786),
    @(0.2333249530443311, 0.1969762600297997, 0.1636653930274292, 1.341072836250262, 0.7382829981567127, 0.8518674884310968, 0.9268513112354491, 0.1806385917172122, 1.577459891748731, 0.2141990775815827),
    @(0.2331273504922962, 0.196770368163186, 0.1636105993333459, 1.341718477567156, 0.7388166395949938, 0.8523417643582221, 0.9273287121718958, 0.1806529980754163, 1.568340707861751, 0.2141835887719665),
    @(0.2345054055106885, 0.1982037370876526, 0.1639958991909687, 1.337310595660554, 0.7351649016993775, 0.8490885603913654, 0.924067057863688, 0.1805595117404089, 1.631595292377199, 0.2142978310937167),
    @(0.2407458061142052, 0.2046332280699659, 0.16581875780874, 1.319680891285145, 0.720346954216204, 0.8356928933832393, 0.9109622178944257, 0.1803078373649569, 1.909642322352568, 0.2149681952694849),
    @(0.2535477889160518, 0.2176084131811962, 0.1698304751258455, 1.291662254932497, 0.6959653707740188, 0.8128759042197089, 0.889899680405577, 0.1803905429211099, 2.450750767957175, 0.2168769773671642),
    @(0.2632917979076979, 0.2273591965795276, 0.1730417211875697, 1.275078940004718, 0.680918426072374, 0.798217524768468, 0.8772572619388583, 0.1808011713499837, 2.845603361750193, 0.2186391306503808),
    @(0.2677969314067639, 0.2318412462856401, 0.1745593646504346, 1.268405280550468, 0.6746988246353567, 0.7920063777566497, 0.8721221873420646, 0.1810636612890235, 3.024611681306851, 0.219518403573673),
    @(0.2695132188805758, 0.2335450308488873, 0.1751421759679381, 1.266003393746885, 0.6724338206400873, 0.7897201272066923, 0.8702663598307652, 0.1811739162004002, 3.092305804920329, 0.219862480498108),
    @(0.2691431301724378, 0.2331778021459883, 0.1750162972567182, 1.266515108619245, 0.6729176122383578, 0.7902095858625771, 0.8706620978612918, 0.1811496884831669, 3.077730845497797, 0.2197878838465357),
    @(0.2679379259542998, 0.2319812875937544, 0.1746071506389271, 1.268205164143879, 0.6745106714419791, 0.7918169679747038, 0.8719677285981646, 0.1810725146284753, 3.030182795438918, 0.2195464885191711),
    @(0.2672010396673841, 0.2312492341394545, 0.1743575912654336, 1.26925669204212, 0.6754982249901147, 0.7928101035884154, 0.8727790221200493, 0.1810266562859653, 3.001046072724648, 0.2194000729131744),
    @(0.2629988245080881, 0.2270672053256817, 0.1729436781115936, 1.275532602605864, 0.6813375007756264, 0.7986326391553575, 0.8776052582916307, 0.1807855376807623, 2.833892050888267, 0.2185832263646361),
    @(0.2604393802764093, 0.2245134438115031, 0.1720908016090164, 1.279605652419029, 0.6850800849506271, 0.8023216711126508, 0.8807238696319999, 0.1806569854434343, 2.73118853664937, 0.218101971043005),
    @(0.2589740912657987, 0.2230489604604315, 0.1716056024062524, 1.282030258750986, 0.6872915613486441, 0.8044865185432002, 0.8825755878720614, 0.1805901719042424, 2.67205883873703, 0.2178324785974155),
    @(0.2584791478608111, 0.2225538668063081, 0.1714422434624545, 1.282865249782496, 0.6880504260068534, 0.805226884624588, 0.8832125012119576, 0.1805687748197045, 2.652028824089371, 0.2177424906725633),
    @(0.2607111310356629, 0.2247848445274769, 0.1721810382065705, 1.279163592287553, 0.6846755884588021, 0.8019245151822147, 0.8803858869386119, 0.180669932711794, 2.742127469058005, 0.2181524450155266),
    @(0.2682916454542976, 0.2323325568076058, 0.1747271072168601, 1.267705352505686, 0.6740403002796924, 0.791343055291641, 0.8715818245037639, 0.1810948880705681, 3.044151354184635, 0.219617090925702),
    @(0.2733058706101019, 0.2373034190599981, 0.1764383698591203, 1.26094704669557, 0.6676155587955819, 0.7848108787374883, 0.8663450121667964, 0.1814358827389668, 3.241000580089462, 0.2206390863870595),
    @(0.2706242473656744, 0.2346469451728836, 0.175520732192151, 1.264487207507251, 0.6709963325893185, 0.788262122751199, 0.8690926353940327, 0.1812481085413395, 3.135989379807597, 0.2200877193873083),
    @(0.2605882532905213, 0.2246621327041112, 0.1721402262547329, 1.279363189282414, 0.6848582750255119, 0.8021039324057213, 0.8805385058059514, 0.1806640571559228, 2.737182239279775, 0.2181296033405218),
    @(0.2500246725670792, 0.2140595976431285, 0.1686986776733015, 1.2985397488901, 0.7020590035659353, 0.8186788634416473, 0.8951007145294376, 0.1803066262632598, 2.304828220775619, 0.2162972439707005)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowVals = $rows[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item(2 + $r, 3 + $c).Value = $rowVals[$c]
    }
}

# Column O for rows 2-25
$oVals = @(
    3.100743789450632,
    3.145108458399278,
    3.174495473503839,
    3.187010183775627,
    3.189120794403607,
    3.174662068608257,
    3.115594775665855,
    3.016825054838364,
    2.954694707323256,
    2.928703460474381,
    2.919188638728613,
    2.921223248066269,
    2.927914102418057,
    2.932055116634572,
    2.956439083630769,
    2.971980358252154,
    2.981133109289885,
    2.984268774089514,
    2.970303827811094,
    2.925939941987352,
    2.898854847197725,
    2.913135721910223,
    2.971061108641095,
    3.041714757884634
)

for ($r = 0; $r -lt $oVals.Length; $r++) {
    $ws.Cells.Item(2 + $r, 15).Value = $oVals[$r]
}
